$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.540.92'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '2.068.73'
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.622'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.61%  '
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.71'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.389'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0781'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.82'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.51%  '
$ws.Range("D13").Value = '2.378.52'
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.10'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.763'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.32'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.50%  '
$ws.Range("D17").Value = '2.077.07'
$ws.Range("E17").Value = '  -0.68%  '
$ws.Range("D18").Value = '37.483.57'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.31%  '
$ws.Range("D21").Value = '0.0₃0828'
$ws.Range("E21").Value = '  -2.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.36'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.17%  '
$ws.Range("E28").Value = '  -4.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("E30").Value = '  -4.62%  '
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.59'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0630'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.87%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.51'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.45%  '
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.31'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.82%  '
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.32'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.14%  '
$ws.Range("E40").Value = '  +3.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.24'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.93%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.92'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.23%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.490.84'
$ws.Range("E43").Value = '  +2.69%  '
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0955'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.99'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.73%  '
$ws.Range("E46").Value = '  +2.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.06'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.27%  '
$ws.Range("E48").Value = '  -2.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.27'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.96'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.21%  '
$ws.Range("D51").Value = '2.264.36'
$ws.Range("E51").Value = '  -0.52%  '
